$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Text columns (A-D). Prefix with a leading apostrophe so Excel stores the
# literal text instead of auto-converting date/time/number-looking strings,
# then reset the style back to the default so no stray "quote prefix" /
# number-format style sticks to the cell.
$ws.Range("A33").Value = "'2023-06-08"
$ws.Range("A33").Style = "常规"

$ws.Range("B33").Value = "'21:13:22"
$ws.Range("B33").Style = "常规"

$ws.Range("C33").Value = "'Thursday"
$ws.Range("C33").Style = "常规"

$ws.Range("D33").Value = "'23"
$ws.Range("D33").Style = "常规"

# Numeric columns (E-T)
$ws.Range("E33").Value = 120091
$ws.Range("F33").Value = 134304
$ws.Range("G33").Value = 160552
$ws.Range("H33").Value = 131525
$ws.Range("I33").Value = 175605
$ws.Range("J33").Value = 113315
$ws.Range("K33").Value = 201357
$ws.Range("L33").Value = 221438
$ws.Range("M33").Value = 172959
$ws.Range("N33").Value = 120171
$ws.Range("O33").Value = 38667
$ws.Range("P33").Value = 34433
$ws.Range("Q33").Value = 50886
$ws.Range("R33").Value = -1
$ws.Range("S33").Value = 36961
$ws.Range("T33").Value = -1
